# "zname chyby v zadavani planu jsou odstraneny"
# Add a "koupeno_old" column (H) capturing the previous purchase-plan
# quantities/flags, and correct the "potreba" (E) quantities that had been
# mis-entered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header H1, matching the style of the other header cells (G1) ---
$ws.Range("H1").Value = "koupeno_old"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# --- Corrected "potreba" values in column E ---
$eValues = @{
    2  = 365.664
    3  = 365.664
    4  = 3.6568
    5  = 10.97
    6  = 17.1864
    7  = 4000
    8  = 266.668
    9  = 266.668
    10 = 4000
    11 = 731.328
    12 = 731.328
    13 = 7.313600000000001
    14 = 21.94
    15 = 34.3728
    16 = 8000
    17 = 533.336
    18 = 533.336
    19 = 8000
}

foreach ($row in $eValues.Keys) {
    $ws.Cells.Item($row, 5).Value = $eValues[$row]
}

# --- Row 4 is the one exception: "koupeno" flips to FALSE and its
#     "koupeno_old" becomes an empty marker (no prior purchase) ---
$ws.Cells.Item(4, 7).Value = $false
$ws.Cells.Item(4, 8).Formula = "=""""" 

# --- Every other row keeps "koupeno" TRUE and now also records
#     "koupeno_old" = TRUE ---
foreach ($row in 2..19) {
    if ($row -eq 4) { continue }
    $ws.Cells.Item($row, 7).Value = $true
    $ws.Cells.Item($row, 8).Value = $true
}
